# Added filtering options for the Component Analysis
# For each data row (row 2..47 on Sheet1), a new "Q0" error value is
# inserted into column B; the previously-existing values (old B..J) shift
# one column to the right (C..K), and - for the rows that were already full
# (B..K) - the oldest trailing value (old K) drops off the end.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newData = New-Object 'object[,]' 46,10

# Row 2
$newData[0,0] = 0.04231489763667018
$newData[0,1] = -1.112737948256028
$newData[0,2] = 0.426895895633562
$newData[0,3] = 0.1572417756599283
$newData[0,4] = 0.7087734795590415
$newData[0,5] = 0.482846718646189
$newData[0,6] = 0.5906624340427387
$newData[0,7] = 0.4942628133784369
$newData[0,8] = 0.9976736988314454
$newData[0,9] = 1.256195640754894

# Row 3
$newData[1,0] = 1.228061995268202
$newData[1,1] = 0.9584078752945684
$newData[1,2] = 1.509939579193682
$newData[1,3] = 1.284012818280829
$newData[1,4] = 1.391828533677379
$newData[1,5] = 1.295428913013077
$newData[1,6] = 1.798839798466086
$newData[1,7] = 2.057361740389535
$newData[1,8] = 0.3333791554112209
$newData[1,9] = 0.9870183717739222

# Row 4
$newData[2,0] = 0.5137754236260815
$newData[2,1] = 0.2878486627132291
$newData[2,2] = 0.3956643781097787
$newData[2,3] = 0.299264757445477
$newData[2,4] = 0.8026756428984855
$newData[2,5] = 1.061197584821935
$newData[2,6] = -0.6627850001563792
$newData[2,7] = -0.00914578379367792
$newData[2,8] = 0.3543981511147533
$newData[2,9] = -0.5463776937648762

# Row 5
$newData[3,0] = 0.4028038717171413
$newData[3,1] = 0.3064042510528396
$newData[3,2] = 0.8098151365058481
$newData[3,3] = 1.068337078429297
$newData[3,4] = -0.6556455065490165
$newData[3,5] = -0.002006290186315263
$newData[3,6] = 0.3615376447221159
$newData[3,7] = -0.5392382001575136
$newData[3,8] = 0.5311373567961393
$newData[3,9] = -0.04592255722032257

# Row 6
$newData[4,0] = 0.5336388157440486
$newData[4,1] = 0.7921607576674977
$newData[4,2] = -0.9318218273108161
$newData[4,3] = -0.2781826109481148
$newData[4,4] = 0.0853613239603164
$newData[4,5] = -0.8154145209193131
$newData[4,6] = 0.2549610360343397
$newData[4,7] = -0.3220988779821221
$newData[4,8] = -0.5277048620751383
$newData[4,9] = 0.4514401196459589

# Row 7
$newData[5,0] = -1.102517691576566
$newData[5,1] = -0.4488784752138651
$newData[5,2] = -0.08533454030543387
$newData[5,3] = -0.9861103851850633
$newData[5,4] = 0.08426517176858947
$newData[5,5] = -0.4927947422478724
$newData[5,6] = -0.6984007263408886
$newData[5,7] = 0.2807442553802086
$newData[5,8] = -0.4556020879571533
$newData[5,9] = -0.4735719026113467

# Row 8
$newData[6,0] = 0.1737007515684039
$newData[6,1] = -0.7270750933112256
$newData[6,2] = 0.3433004636424272
$newData[6,3] = -0.2337594503740346
$newData[6,4] = -0.4393654344670508
$newData[6,5] = 0.5397795472540464
$newData[6,6] = -0.1965667960833156
$newData[6,7] = -0.2145366107375089
$newData[6,8] = -0.1911317451971721
$newData[6,9] = -0.5563061809639129

# Row 9
$newData[7,0] = 0.8389029408811082
$newData[7,1] = 0.2618430268646463
$newData[7,2] = 0.05623704277163011
$newData[7,3] = 1.035382024492727
$newData[7,4] = 0.2990356811553654
$newData[7,5] = 0.281065866501172
$newData[7,6] = 0.3044707320415089
$newData[7,7] = -0.06070370372523193
$newData[7,8] = 0.5951158134457202
$newData[7,9] = 0.4642449325042965

# Row 10
$newData[8,0] = -0.6954484448595206
$newData[8,1] = 0.2836965368615766
$newData[8,2] = -0.4526498064757853
$newData[8,3] = -0.4706196211299787
$newData[8,4] = -0.4472147555896419
$newData[8,5] = -0.8123891913563827
$newData[8,6] = -0.1565696741854305
$newData[8,7] = -0.2874405551268542
$newData[8,8] = -0.3969583540509761
$newData[8,9] = -1.226736141118529

# Row 11
$newData[9,0] = -0.2465870357053012
$newData[9,1] = -0.2645568503594945
$newData[9,2] = -0.2411519848191577
$newData[9,3] = -0.6063264205858985
$newData[9,4] = 0.04949309658505369
$newData[9,5] = -0.08137778435637
$newData[9,6] = -0.1908955832804919
$newData[9,7] = -1.020673370348045
$newData[9,8] = -0.4401475811793975
$newData[9,9] = -0.5235336826091774

# Row 12
$newData[10,0] = -0.1975260465718366
$newData[10,1] = -0.5627004823385774
$newData[10,2] = 0.09311903483237477
$newData[10,3] = -0.03775184610904891
$newData[10,4] = -0.1472696450331708
$newData[10,5] = -0.9770474321007241
$newData[10,6] = -0.3965216429320764
$newData[10,7] = -0.4799077443618563
$newData[10,8] = -0.697817621901615
$newData[10,9] = 0.8659601817345554

# Row 13
$newData[11,0] = 0.4425040297996861
$newData[11,1] = 0.3116331488582624
$newData[11,2] = 0.2021153499341405
$newData[11,3] = -0.6276624371334127
$newData[11,4] = -0.04713664796476502
$newData[11,5] = -0.130522749394545
$newData[11,6] = -0.3484326269343037
$newData[11,7] = 1.215345176701867
$newData[11,8] = -0.1206385293801969
$newData[11,9] = 1.319975887149931

# Row 14
$newData[12,0] = -0.2720610750631522
$newData[12,1] = -1.101838862130705
$newData[12,2] = -0.5213130729620578
$newData[12,3] = -0.6046991743918377
$newData[12,4] = -0.8226090519315964
$newData[12,5] = 0.741168751704574
$newData[12,6] = -0.5948149543774897
$newData[12,7] = 0.8457994621526386
$newData[12,8] = -0.2546446521081841
$newData[12,9] = -0.464919050277854

# Row 15
$newData[13,0] = -0.1065518669046048
$newData[13,1] = -0.1899379683343848
$newData[13,2] = -0.4078478458741435
$newData[13,3] = 1.155929957762027
$newData[13,4] = -0.1800537483200367
$newData[13,5] = 1.260560668210092
$newData[13,6] = 0.1601165539492688
$newData[13,7] = -0.05015784422040104
$newData[13,8] = 0.5483863591633293
$newData[13,9] = 0.2227177183881149

# Row 16
$newData[14,0] = -0.1895682054566924
$newData[14,1] = 1.374209598179478
$newData[14,2] = 0.03822589209741434
$newData[14,3] = 1.478840308627543
$newData[14,4] = 0.3783961943667199
$newData[14,5] = 0.16812179619705
$newData[14,6] = 0.7666659995807804
$newData[14,7] = 0.440997358805566
$newData[14,8] = 1.15817716018117
$newData[14,9] = 2.854049435818871

# Row 17
$newData[15,0] = 1.157000698704573
$newData[15,1] = -0.1789830073774904
$newData[15,2] = 1.261631409152638
$newData[15,3] = 0.1611872948918152
$newData[15,4] = -0.0490871032778547
$newData[15,5] = 0.5494571001058757
$newData[15,6] = 0.2237884593306613
$newData[15,7] = 0.9409682607062649
$newData[15,8] = 2.636840536343966
$newData[15,9] = 9.481087331268544

# Row 18
$newData[16,0] = -0.4886691766355519
$newData[16,1] = 0.9519452398945764
$newData[16,2] = -0.1484988743662463
$newData[16,3] = -0.3587732725359162
$newData[16,4] = 0.2397709308478142
$newData[16,5] = -0.0858977099274002
$newData[16,6] = 0.6312820914482035
$newData[16,7] = 2.327154367085904
$newData[16,8] = 9.171401162010483
$newData[16,9] = -8.379252921090828

# Row 19
$newData[17,0] = 1.10624937372658
$newData[17,1] = 0.005805259465757717
$newData[17,2] = -0.2044691387039121
$newData[17,3] = 0.3940750646798182
$newData[17,4] = 0.06840642390460383
$newData[17,5] = 0.7855862252802075
$newData[17,6] = 2.481458500917908
$newData[17,7] = 9.325705295842486
$newData[17,8] = -8.224948787258825
$newData[17,9] = -0.5155521876026774

# Row 20
$newData[18,0] = -0.6446211617534254
$newData[18,1] = -0.8548955599230954
$newData[18,2] = -0.256351356539365
$newData[18,3] = -0.5820199973145794
$newData[18,4] = 0.1351598040610243
$newData[18,5] = 1.831032079698725
$newData[18,6] = 8.675278874623302
$newData[18,7] = -8.875375208478008
$newData[18,8] = -1.165978608821861
$newData[18,9] = 0.416934052065051

# Row 21
$newData[19,0] = -0.6387305113048862
$newData[19,1] = -0.04018630792115581
$newData[19,2] = -0.3658549486963703
$newData[19,3] = 0.3513248526792334
$newData[19,4] = 2.047197128316934
$newData[19,5] = 8.891443923241512
$newData[19,6] = -8.659210159859798
$newData[19,7] = -0.9498135602036515
$newData[19,8] = 0.6330991006832601
$newData[19,9] = -2.333187891665996

# Row 22
$newData[20,0] = 0.3668428211138005
$newData[20,1] = 0.04117418033858611
$newData[20,2] = 0.7583539817141898
$newData[20,3] = 2.454226257351891
$newData[20,4] = 9.298473052276469
$newData[20,5] = -8.252181030824842
$newData[20,6] = -0.5427844311686951
$newData[20,7] = 1.040128229718217
$newData[20,8] = -1.926158762631039
$newData[20,9] = 0.3378723920334236

# Row 23
$newData[21,0] = -0.4578680368388337
$newData[21,1] = 0.25931176453677
$newData[21,2] = 1.955184040174471
$newData[21,3] = 8.799430835099049
$newData[21,4] = -8.751223248002262
$newData[21,5] = -1.041826648346115
$newData[21,6] = 0.5410860125407967
$newData[21,7] = -2.425200979808459
$newData[21,8] = -0.1611698251439962
$newData[21,9] = -0.6123239800707349

# Row 24
$newData[22,0] = 0.4181606776922825
$newData[22,1] = 2.114032953329983
$newData[22,2] = 8.958279748254562
$newData[22,3] = -8.592374334846749
$newData[22,4] = -0.8829777351906024
$newData[22,5] = 0.6999349256963092
$newData[22,6] = -2.266352066652947
$newData[22,7] = -0.002320911988483623
$newData[22,8] = -0.4534750669152223
$newData[22,9] = -0.5938391304291964

# Row 25
$newData[23,0] = 2.057869132359739
$newData[23,1] = 8.902115927284317
$newData[23,2] = -8.648538155816993
$newData[23,3] = -0.9391415561608464
$newData[23,4] = 0.6437711047260652
$newData[23,5] = -2.322515887623191
$newData[23,6] = -0.05848473295872768
$newData[23,7] = -0.5096388878854663
$newData[23,8] = -0.6500029513994404
$newData[23,9] = -0.1237639405372229

# Row 26
$newData[24,0] = 6.652313087672924
$newData[24,1] = -10.89834099542839
$newData[24,2] = -3.188944395772239
$newData[24,3] = -1.606031734885327
$newData[24,4] = -4.572318727234583
$newData[24,5] = -2.30828757257012
$newData[24,6] = -2.759441727496859
$newData[24,7] = -2.899805791010833
$newData[24,8] = -2.373566780148615
$newData[24,9] = -2.515395726887889

# Row 27
$newData[25,0] = -18.36749132628568
$newData[25,1] = -10.65809472662953
$newData[25,2] = -9.07518206574262
$newData[25,3] = -12.04146905809188
$newData[25,4] = -9.777437903427414
$newData[25,5] = -10.22859205835415
$newData[25,6] = -10.36895612186813
$newData[25,7] = -9.842717111005909
$newData[25,8] = -9.984546057745183
$newData[25,9] = -9.341571436343118

# Row 28
$newData[26,0] = 7.513167073507937
$newData[26,1] = 9.096079734394849
$newData[26,2] = 6.129792742045593
$newData[26,3] = 8.393823896710055
$newData[26,4] = 7.942669741783317
$newData[26,5] = 7.802305678269343
$newData[26,6] = 8.32854468913156
$newData[26,7] = 8.186715742392286
$newData[26,8] = 8.82969036379435
$newData[26,9] = 8.954103070123933

# Row 29
$newData[27,0] = 0.9564081874156993
$newData[27,1] = -2.009878804933557
$newData[27,2] = 0.2541523497309064
$newData[27,3] = -0.1970018051958322
$newData[27,4] = -0.3373658687098063
$newData[27,5] = 0.1888731421524112
$newData[27,6] = 0.0470441954131372
$newData[27,7] = 0.690018816815201
$newData[27,8] = 0.8144315231447835
$newData[27,9] = 0.4168013171071663

# Row 30
$newData[28,0] = -4.157449276732949
$newData[28,1] = -1.893418122068486
$newData[28,2] = -2.344572276995224
$newData[28,3] = -2.484936340509199
$newData[28,4] = -1.958697329646981
$newData[28,5] = -2.100526276386255
$newData[28,6] = -1.457551654984191
$newData[28,7] = -1.333138948654609
$newData[28,8] = -1.730769154692226
$newData[28,9] = -1.806937868393072

# Row 31
$newData[29,0] = 1.546611864454844
$newData[29,1] = 1.095457709528105
$newData[29,2] = 0.9550936460141312
$newData[29,3] = 1.481332656876349
$newData[29,4] = 1.339503710137075
$newData[29,5] = 1.982478331539139
$newData[29,6] = 2.106891037868721
$newData[29,7] = 1.709260831831104
$newData[29,8] = 1.633092118130258
$newData[29,9] = 1.90941417625318

# Row 32
$newData[30,0] = 1.156631887942306
$newData[30,1] = 1.016267824428332
$newData[30,2] = 1.54250683529055
$newData[30,3] = 1.400677888551276
$newData[30,4] = 2.04365250995334
$newData[30,5] = 2.168065216282922
$newData[30,6] = 1.770435010245305
$newData[30,7] = 1.694266296544459
$newData[30,8] = 1.970588354667381
$newData[30,9] = 1.799368532355657

# Row 33
$newData[31,0] = -1.025188112727922
$newData[31,1] = -0.4989491018657047
$newData[31,2] = -0.6407780486049788
$newData[31,3] = 0.002196572797085183
$newData[31,4] = 0.1266092791266676
$newData[31,5] = -0.2710209269109496
$newData[31,6] = -0.3471896406117954
$newData[31,7] = -0.07086758248887381
$newData[31,8] = -0.2420874048005978
$newData[31,9] = -0.08893999929185659

# Row 34
$newData[32,0] = 0.0836454351679363
$newData[32,1] = -0.05818351157133772
$newData[32,2] = 0.5847911098307261
$newData[32,3] = 0.7092038161603086
$newData[32,4] = 0.3115736101226914
$newData[32,5] = 0.2354048964218456
$newData[32,6] = 0.5117269545447671
$newData[32,7] = 0.3405071322330432
$newData[32,8] = 0.4936545377417844
$newData[32,9] = 0.2162257633186657

# Row 35
$newData[33,0] = -0.1538585523806955
$newData[33,1] = 0.4891160690213684
$newData[33,2] = 0.6135287753509509
$newData[33,3] = 0.2158985693133336
$newData[33,4] = 0.1397298556124878
$newData[33,5] = 0.4160519137354094
$newData[33,6] = 0.2448320914236854
$newData[33,7] = 0.3979794969324266
$newData[33,8] = 0.1205507225093079
$newData[33,9] = -0.04243697084963852

# Row 36
$newData[34,0] = 0.7495351060200912
$newData[34,1] = 0.8739478123496736
$newData[34,2] = 0.4763176063120564
$newData[34,3] = 0.4001488926112106
$newData[34,4] = 0.6764709507341322
$newData[34,5] = 0.5052511284224082
$newData[34,6] = 0.6583985339311494
$newData[34,7] = 0.3809697595080307
$newData[34,8] = 0.2179820661490843
$newData[34,9] = 0.0938994907545665

# Row 37
$newData[35,0] = 0.03849281619118239
$newData[35,1] = -0.3591373898464348
$newData[35,2] = -0.4353061035472806
$newData[35,3] = -0.158984045424359
$newData[35,4] = -0.330203867736083
$newData[35,5] = -0.1770564622273418
$newData[35,6] = -0.4544852366504605
$newData[35,7] = -0.617472930009407
$newData[35,8] = -0.7415555054039247
$newData[35,9] = -0.2254024683979639

# Row 38
$newData[36,0] = -0.2590580299438133
$newData[36,1] = -0.3352267436446591
$newData[36,2] = -0.0589046855217375
$newData[36,3] = -0.2301245078334615
$newData[36,4] = -0.07697710232472027
$newData[36,5] = -0.354405876747839
$newData[36,6] = -0.5173935701067854
$newData[36,7] = -0.6414761455013032
$newData[36,8] = -0.1253231084953424
$newData[36,9] = -0.3352267436446591

# Row 39
$newData[37,0] = 0.01855976243503714
$newData[37,1] = 0.2948818205579588
$newData[37,2] = 0.1236619982462347
$newData[37,3] = 0.276809403754976
$newData[37,4] = -0.0006193706681427817
$newData[37,5] = -0.1636070640270892
$newData[37,6] = -0.287689639421607
$newData[37,7] = 0.2284633975843539
$newData[37,8] = 0.01855976243503714
$newData[37,9] = $null

# Row 40
$newData[38,0] = 0.1467044301255134
$newData[38,1] = -0.0245153921862106
$newData[38,2] = 0.1286320133225306
$newData[38,3] = -0.1487967611005881
$newData[38,4] = -0.3117844544595345
$newData[38,5] = -0.4358670298540523
$newData[38,6] = 0.0802860071519085
$newData[38,7] = -0.1296176279974082
$newData[38,8] = $null
$newData[38,9] = $null

# Row 41
$newData[39,0] = -0.1819613811903656
$newData[39,1] = -0.02881397568162436
$newData[39,2] = -0.3062427501047431
$newData[39,3] = -0.4692304434636895
$newData[39,4] = -0.5933130188582073
$newData[39,5] = -0.07715998185224648
$newData[39,6] = -0.2870636170015632
$newData[39,7] = $null
$newData[39,8] = $null
$newData[39,9] = $null

# Row 42
$newData[40,0] = 0.4718454808444464
$newData[40,1] = 0.1944167064213277
$newData[40,2] = 0.0314290130623813
$newData[40,3] = -0.09265356233213651
$newData[40,4] = 0.4234994746738243
$newData[40,5] = 0.2135958395245076
$newData[40,6] = $null
$newData[40,7] = $null
$newData[40,8] = $null
$newData[40,9] = $null

# Row 43
$newData[41,0] = -0.08594117411414147
$newData[41,1] = -0.2489288674730878
$newData[41,2] = -0.3730114428676057
$newData[41,3] = 0.1431415941383551
$newData[41,4] = -0.06676204101096155
$newData[41,5] = $null
$newData[41,6] = $null
$newData[41,7] = $null
$newData[41,8] = $null
$newData[41,9] = $null

# Row 44
$newData[42,0] = -0.07695400962807622
$newData[42,1] = -0.201036585022594
$newData[42,2] = 0.3151164519833668
$newData[42,3] = 0.1052128168340501
$newData[42,4] = $null
$newData[42,5] = $null
$newData[42,6] = $null
$newData[42,7] = $null
$newData[42,8] = $null
$newData[42,9] = $null

# Row 45
$newData[43,0] = -0.5068991247689255
$newData[43,1] = 0.009253912237035311
$newData[43,2] = -0.2006497229122814
$newData[43,3] = $null
$newData[43,4] = $null
$newData[43,5] = $null
$newData[43,6] = $null
$newData[43,7] = $null
$newData[43,8] = $null
$newData[43,9] = $null

# Row 46
$newData[44,0] = 0.6215838649243215
$newData[44,1] = 0.4116802297750048
$newData[44,2] = $null
$newData[44,3] = $null
$newData[44,4] = $null
$newData[44,5] = $null
$newData[44,6] = $null
$newData[44,7] = $null
$newData[44,8] = $null
$newData[44,9] = $null

# Row 47
$newData[45,0] = -0.2766911554241067
$newData[45,1] = $null
$newData[45,2] = $null
$newData[45,3] = $null
$newData[45,4] = $null
$newData[45,5] = $null
$newData[45,6] = $null
$newData[45,7] = $null
$newData[45,8] = $null
$newData[45,9] = $null

$ws.Range("B2:K47").Value = $newData
